$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22 (pushes existing rows 22-55 down to 23-56,
# carrying their values/styles with them, exactly like Excel's UI "Insert Row").
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new price-report record.
$ws.Cells.Item(22, 1).Value = 10
$ws.Cells.Item(22, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(22, 3).Value = "La Araucanía"
$ws.Cells.Item(22, 4).Value = 44895
$ws.Cells.Item(22, 5).Value = 9
$ws.Cells.Item(22, 6).Value = 100112042
$ws.Cells.Item(22, 7).Value = "Locoto"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 40
$ws.Cells.Item(22, 11).Value = 2500
$ws.Cells.Item(22, 12).Value = 2500
$ws.Cells.Item(22, 13).Value = 2500
$ws.Cells.Item(22, 14).Value = '$/kilo'
$ws.Cells.Item(22, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(22, 16).Value = 2500
$ws.Cells.Item(22, 17).Value = 1
$ws.Cells.Item(22, 18).Value = "Hortaliza"
